# Update cryptos list values per commit "Updated cryptos list on Wed Sep 27 14:14:16 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-47: update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Price values that look like plain numbers must be forced back to text (matching the
# original "General"-styled inline-string cells) so Excel does not auto-convert them to numbers.
$ws.Range("D2").Value = '26.627.76'
$ws.Range("E2").Value = '  +1.66%  '
$ws.Range("D3").Value = '1.620.27'
$ws.Range("E3").Value = '  +2.11%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = '  +0.70%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  +0.58%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.28'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.03%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0858'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.23%  '
$ws.Range("D12").Value = '1.848.71'
$ws.Range("E12").Value = '  +2.14%  '
$ws.Range("D13").Value = '1.606.86'
$ws.Range("E13").Value = '  +1.55%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.04'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.77%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '64.54'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.71%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.510'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.92%  '
$ws.Range("D17").Value = '26.641.62'
$ws.Range("E17").Value = '  +1.65%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '232.68'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +9.11%  '
$ws.Range("D19").Value = '0.0₃0730'
$ws.Range("E19").Value = '  +0.52%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.69'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.73%  '
$ws.Range("E21").Value = '  -0.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.38'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.23'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.23%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.07'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.67'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.47%  '
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.04'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.57%  '
$ws.Range("E28").Value = '  +2.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.62'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.91%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0497'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.21%  '
$ws.Range("E31").Value = '  +1.51%  '
$ws.Range("E32").Value = '  +1.76%  '
$ws.Range("D33").Value = '1.457.09'
$ws.Range("E33").Value = '  +8.77%  '
$ws.Range("E34").Value = '  +1.70%  '
$ws.Range("E35").Value = '  -0.53%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.48'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.78%  '
$ws.Range("E37").Value = '  -1.96%  '
$ws.Range("E38").Value = '  +0.00%  '
$ws.Range("E39").Value = '  +2.00%  '
$ws.Range("E40").Value = '  +2.41%  '
$ws.Range("E41").Value = '  +0.07%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.20'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.02%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.948'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.77%  '
$ws.Range("D44").Value = '1.759.46'
$ws.Range("E44").Value = '  +2.15%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.766'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.40%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '61.72'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.99%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '88.34'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.07%  '

# Rows 48-51 shift up: the BabyDogeCoin row was dropped and EnergySwap was appended at the
# bottom, so RenderToken/Cronos/Algorand data each moved up one row and row 51 becomes EnergySwap.
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.50'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.38%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0504'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.61%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0962'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.94%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.49'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.31%  '
